# Commit: "Added Apache PO for Excel Reading, and SearchTest Case"
# Rename Sheet1 -> ableToSearchProducts and populate the small product table
# used by the new search test.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "ableToSearchProducts"

# Row 1: headers
$ws.Range("A1").Value = "Product Name"
$ws.Range("B1").Value = "Product  Price"

# Row 2
$ws.Range("A2").Value = "Canon EOS 5D"
# The price looks like a currency number ("$98.00"), so Excel's normal
# literal-entry parsing would store it as a numeric/currency cell. Flip the
# cell to Text just long enough to type the literal string, then flip back
# to General so the cell keeps using the sheet's default (General) style -
# this matches how the source workbook stores these as plain text cells.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "$98.00"
$ws.Range("B2").NumberFormat = "General"

# Row 3
$ws.Range("A3").Value = "Samsung Galaxy Tab 10.1"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "$241.99"
$ws.Range("B3").NumberFormat = "General"

# Row 4
$ws.Range("A4").Value = "MacBook Air"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "$1,202.00"
$ws.Range("B4").NumberFormat = "General"
